$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "lrc-panel-result" row (row 5) belongs right after the header, ahead of
# the "lrc-result" group (rows 2-4). Fix the ordering by inserting a blank
# row at row 2 (pushing the existing rows 2-4 down to 3-5), copying row 5's
# (now row 6's) content into the new row 2, then removing the now-duplicate
# trailing row.
$ws.Range("A2:K2").Insert()
$ws.Range("A6:K6").Copy($ws.Range("A2"))
$ws.Range("A6:K6").Delete()
